# "Added Sprint 2 docs" — roll the Sprint1 burn-down tracker over to Sprint2:
# bump the starting backlog total and fill in the daily
# completed/backlog-change numbers that were captured during the sprint.
# The dependent BurnDn/Ideal formulas in columns E/F (and the chart that
# plots them) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New sprint starting backlog (Backlog column, row 2)
$ws.Range("C2").Value = 45

# Daily Backlog-added (C) / Completed (D) entries recorded through the sprint
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1

$ws.Range("D17").Value = 8

$ws.Range("C19").Value = -1
$ws.Range("D19").Value = 8

$ws.Range("D20").Value = 14

$ws.Range("D21").Value = 6

$ws.Range("D22").Value = 5

# Chart title: Sprint1 -> Sprint2
$co = $ws.ChartObjects("Chart 3")
$co.Chart.ChartTitle.Text = "BURN DOWN CHART      6733 Team2 Sprint2"

# Leave the selection where the author last left it
$ws.Range("E29").Select() | Out-Null
